# Realestate Update resale numbers 2024-01-12 22:59
# Append the new resale-number row (row 53) for 2024-01-12 22:59:29.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 53

# Columns A, B and D hold date-/time-/week-looking text ("2024-01-12",
# "22:59:29", "01") that Excel's input parser would otherwise coerce into
# a date serial / time serial / number. Pre-format those cells as Text so
# the literal strings are preserved, matching the source data (which
# stores them as plain text). Column C ("Friday") parses fine as text on
# its own and needs no special handling.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("B$row").NumberFormat = "@"
$ws.Range("D$row").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-12"
$ws.Cells.Item($row, 2).Value = "22:59:29"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "01"

# Columns E-T are the numeric city resale-number values.
$ws.Cells.Item($row, 5).Value  = 136724
$ws.Cells.Item($row, 6).Value  = 142770
$ws.Cells.Item($row, 7).Value  = 172094
$ws.Cells.Item($row, 8).Value  = 148302
$ws.Cells.Item($row, 9).Value  = -1
$ws.Cells.Item($row, 10).Value = 119780
$ws.Cells.Item($row, 11).Value = 225050
$ws.Cells.Item($row, 12).Value = 253075
$ws.Cells.Item($row, 13).Value = 184958
$ws.Cells.Item($row, 14).Value = 110466
$ws.Cells.Item($row, 15).Value = 40950
$ws.Cells.Item($row, 16).Value = 30924
$ws.Cells.Item($row, 17).Value = 73113
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42815
$ws.Cells.Item($row, 20).Value = -1
